# "add tag to example"
#
# Insert a new row above the current header row (row 7) on the active sheet
# and label it "tag" in column A. Excel's normal EntireRow.Insert() behavior
# shifts everything from the old row 7 down (header row -> row 8, all the
# data rows -> one row further down, and the sheet's used range/dimension
# grows by one row accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7, pushing the header + data rows down by one.
$ws.Rows.Item(7).Insert()

# Label the new row.
$ws.Range("A7").Value = "tag"

# Match the row height used throughout the sheet for data rows.
$ws.Rows.Item(7).RowHeight = 15.75

# Leave the selection where the author's last save left it.
$ws.Range("A8").Select()
